$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value2 = 0.1641791044776119
$ws.Range("C2").Value2 = 0.6
$ws.Range("J2").Value2 = 0.01492537313432836
$ws.Range("P2").Value2 = 0.1194029850746269
$ws.Range("S2").Value2 = 0.1014925373134328
$ws.Range("B3").Value2 = 0.01913875598086124
$ws.Range("C3").Value2 = 0.04784688995215311
$ws.Range("J3").Value2 = 0.02870813397129187
$ws.Range("P3").Value2 = 0.7129186602870813
$ws.Range("S3").Value2 = 0.1913875598086124
$ws.Range("J4").Value2 = 0.02173913043478261
$ws.Range("P4").Value2 = 0.7173913043478261
$ws.Range("S4").Value2 = 0.2608695652173913
$ws.Range("B6").Value2 = 0.07224334600760456
$ws.Range("D6").Value2 = 0.007604562737642586
$ws.Range("F6").Value2 = 0.06463878326996197
$ws.Range("J6").Value2 = 0.285171102661597
$ws.Range("O6").Value2 = 0.03422053231939164
$ws.Range("Q6").Value2 = 0.1673003802281369
$ws.Range("R6").Value2 = 0.0532319391634981
$ws.Range("S6").Value2 = 0.3155893536121673
$ws.Range("B7").Value2 = 0.09716599190283401
$ws.Range("D7").Value2 = 0.01619433198380567
$ws.Range("F7").Value2 = 0.08502024291497975
$ws.Range("J7").Value2 = 0.08097165991902834
$ws.Range("O7").Value2 = 0.01619433198380567
$ws.Range("Q7").Value2 = 0.1902834008097166
$ws.Range("R7").Value2 = 0.06072874493927125
$ws.Range("S7").Value2 = 0.4534412955465587
$ws.Range("B8").Value2 = 0.1116504854368932
$ws.Range("D8").Value2 = 0.01456310679611651
$ws.Range("E8").Value2 = 0.002427184466019417
$ws.Range("F8").Value2 = 0.04854368932038835
$ws.Range("J8").Value2 = 0.1359223300970874
$ws.Range("O8").Value2 = 0.01213592233009709
$ws.Range("Q8").Value2 = 0.1699029126213592
$ws.Range("R8").Value2 = 0.1043689320388349
$ws.Range("S8").Value2 = 0.4004854368932039
$ws.Range("B9").Value2 = 0.09012875536480687
$ws.Range("D9").Value2 = 0.004291845493562232
$ws.Range("F9").Value2 = 0.1072961373390558
$ws.Range("J9").Value2 = 0.1373390557939914
$ws.Range("O9").Value2 = 0.008583690987124463
$ws.Range("Q9").Value2 = 0.167381974248927
$ws.Range("R9").Value2 = 0.09012875536480687
$ws.Range("S9").Value2 = 0.3948497854077253
$ws.Range("B10").Value2 = 0.1108144192256342
$ws.Range("D10").Value2 = 0.02336448598130841
$ws.Range("E10").Value2 = 0.001335113484646195
$ws.Range("F10").Value2 = 0.06809078771695594
$ws.Range("J10").Value2 = 0.1388518024032043
$ws.Range("O10").Value2 = 0.01468624833110814
$ws.Range("Q10").Value2 = 0.2242990654205607
$ws.Range("R10").Value2 = 0.08277703604806408
$ws.Range("S10").Value2 = 0.335781041388518
$ws.Range("G11").Value2 = 0.1428571428571428
$ws.Range("J11").Value2 = 0.07908163265306123
$ws.Range("K11").Value2 = 0.2040816326530612
$ws.Range("L11").Value2 = 0.5535714285714286
$ws.Range("S11").Value2 = 0.02040816326530612
$ws.Range("G12").Value2 = 0.7465437788018433
$ws.Range("J12").Value2 = 0.2027649769585254
$ws.Range("K12").Value2 = 0.009216589861751152
$ws.Range("L12").Value2 = 0.01382488479262673
$ws.Range("S12").Value2 = 0.02764976958525346
$ws.Range("G13").Value2 = 0.6153846153846154
$ws.Range("J13").Value2 = 0.3269230769230769
$ws.Range("S13").Value2 = 0.0576923076923077
$ws.Range("F15").Value2 = 0.03947368421052631
$ws.Range("H15").Value2 = 0.1403508771929824
$ws.Range("I15").Value2 = 0.09210526315789473
$ws.Range("J15").Value2 = 0.3596491228070176
$ws.Range("K15").Value2 = 0.04824561403508772
$ws.Range("M15").Value2 = 0.02192982456140351
$ws.Range("O15").Value2 = 0.04824561403508772
$ws.Range("S15").Value2 = 0.25
$ws.Range("F16").Value2 = 0.03720930232558139
$ws.Range("H16").Value2 = 0.1302325581395349
$ws.Range("I16").Value2 = 0.09767441860465116
$ws.Range("J16").Value2 = 0.413953488372093
$ws.Range("K16").Value2 = 0.1116279069767442
$ws.Range("M16").Value2 = 0.03720930232558139
$ws.Range("O16").Value2 = 0.0186046511627907
$ws.Range("S16").Value2 = 0.1534883720930233
$ws.Range("F17").Value2 = 0.01119402985074627
$ws.Range("H17").Value2 = 0.1399253731343283
$ws.Range("I17").Value2 = 0.09328358208955224
$ws.Range("J17").Value2 = 0.4309701492537313
$ws.Range("K17").Value2 = 0.1100746268656716
$ws.Range("M17").Value2 = 0.01865671641791045
$ws.Range("N17").Value2 = 0.001865671641791045
$ws.Range("O17").Value2 = 0.08022388059701492
$ws.Range("S17").Value2 = 0.1138059701492537
$ws.Range("F18").Value2 = 0.01869158878504673
$ws.Range("H18").Value2 = 0.1074766355140187
$ws.Range("I18").Value2 = 0.09813084112149532
$ws.Range("J18").Value2 = 0.4392523364485981
$ws.Range("K18").Value2 = 0.1588785046728972
$ws.Range("M18").Value2 = 0.004672897196261682
$ws.Range("N18").Value2 = 0.004672897196261682
$ws.Range("O18").Value2 = 0.07009345794392523
$ws.Range("S18").Value2 = 0.09813084112149532
$ws.Range("F19").Value2 = 0.01447178002894356
$ws.Range("H19").Value2 = 0.1881331403762663
$ws.Range("I19").Value2 = 0.08827785817655572
$ws.Range("J19").Value2 = 0.3784370477568741
$ws.Range("K19").Value2 = 0.1266280752532561
$ws.Range("M19").Value2 = 0.02170767004341534
$ws.Range("O19").Value2 = 0.06222865412445731
$ws.Range("S19").Value2 = 0.1201157742402316
